$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "mobile"
$ws.Range("B3").Value = "pen"

$ws.Range("E14").Select()
